# Auto-committed on 2023/09/15 週五 17:07:32.90
# Change the column type ("形態") for the CreateDate and LastUpdate rows
# on the "DBD" sheet from DATE to TIMESTAMP, and leave the workbook
# focused/scrolled on that sheet (matching the author's on-save cursor
# position) instead of the "DBS" sheet.

$wb = $excel.ActiveWorkbook

$dbd = $wb.Worksheets.Item("DBD")
$dbs = $wb.Worksheets.Item("DBS")

# CreateDate (row 13) and LastUpdate (row 15) switch from DATE to TIMESTAMP.
$dbd.Range("D13").Value = "TIMESTAMP"
$dbd.Range("D15").Value = "TIMESTAMP"

# Make "DBD" the active sheet/tab (it was "DBS" before) and leave the
# cursor on the last-edited cell, scrolled so row 13 is visible.
$dbd.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$dbd.Range("D15").Select()
